# Scheduled-runner price refresh for Hades_Profits.
#
# The sheet has no formulas -- each leve row's market-price columns
# (H currentAveragePrice / I currentAveragePriceNQ / J currentAveragePriceHQ /
#  K LevePriceNQ / L LevePriceHQ / M LeveProfitNQ / N LeveProfitHQ) are plain
# numeric literals written by an external Universalis price pull, so this
# script just pokes the refreshed numbers straight into the affected cells
# (occasionally a LeveProfit cell goes from 0-profit/blank to populated, or
# vice versa, so a couple of cells are added/cleared instead of just updated).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: "Stuck in the Moment" (Horn Glue)
$ws.Range("H40").Value = 1000.1429
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1000.1429
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1000.1429
$ws.Range("N40").Value = -1350.1429
$ws.Range("M40").ClearContents()

# Row 62: "The Mustache Suits Him" (Enchanted Mythrite Ink)
$ws.Range("H62").Value = 1924.625
$ws.Range("I62").Value = 1970.5714
$ws.Range("K62").Value = 1970.5714
$ws.Range("M62").Value = -1346.5714

# Row 64: "Forged from the Void" (Void Glue)
$ws.Range("H64").Value = 4374.364
$ws.Range("I64").Value = 3433.3333
$ws.Range("J64").Value = 4727.25
$ws.Range("K64").Value = 3433.3333
$ws.Range("L64").Value = 4727.25
$ws.Range("M64").Value = -3185.3333
$ws.Range("N64").Value = -5223.25

# Row 65: "Forgery of Convenience (L)" (Enchanted Mythrite Ink)
$ws.Range("H65").Value = 1924.625
$ws.Range("I65").Value = 1970.5714
$ws.Range("K65").Value = 9852.857
$ws.Range("M65").Value = -6732.857

# Row 67: "Dodging the Draft (L)" (Void Glue)
$ws.Range("H67").Value = 4374.364
$ws.Range("I67").Value = 3433.3333
$ws.Range("J67").Value = 4727.25
$ws.Range("K67").Value = 3433.3333
$ws.Range("L67").Value = 4727.25
$ws.Range("M67").Value = -2575.3333
$ws.Range("N67").Value = -6443.25

# Row 76: "Warding Off Temptation" (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 3283.8096
$ws.Range("I76").Value = 3258.4707
$ws.Range("K76").Value = 3258.4707
$ws.Range("M76").Value = -2943.4707

# Row 79: "The Garden of Arcane Delights (L)" (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 3283.8096
$ws.Range("I79").Value = 3258.4707
$ws.Range("K79").Value = 3258.4707
$ws.Range("M79").Value = -2166.4707

# Row 100: "Asking for a Friend" (Beetle Glue)
$ws.Range("H100").Value = 2988.6365
$ws.Range("I100").Value = 2862.5
$ws.Range("J100").Value = 3140
$ws.Range("K100").Value = 2862.5
$ws.Range("L100").Value = 3140
$ws.Range("M100").Value = -2321.5
$ws.Range("N100").Value = -4222

# Row 128: "Nearly There" (Kumbhiraskin Grimoire)
$ws.Range("H128").Value = 44750
$ws.Range("J128").Value = 44750
$ws.Range("L128").Value = 44750
$ws.Range("N128").Value = -54710

# Row 138: "All-night Crafting" (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 3705684.2
$ws.Range("I138").Value = 1557.25
$ws.Range("J138").Value = 5884582.5
$ws.Range("K138").Value = 4671.75
$ws.Range("L138").Value = 17653747.5
$ws.Range("M138").Value = 468.25
$ws.Range("N138").Value = -17664027.5

$ws = $wb.Worksheets.Item("ARM")
# Row 43: "They've Got Legs" (Steel Sabatons)
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 62: "Hauberk and No Play" (Mythrite Hauberk of Maiming)
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65: "Knights without Armor (L)" (Mythrite Hauberk of Maiming)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 88: "The Mast Chance" (Adamantite Rivets)
$ws.Range("H88").Value = 3267.8823
$ws.Range("I88").Value = 2352
$ws.Range("J88").Value = 3767.4546
$ws.Range("K88").Value = 2352
$ws.Range("L88").Value = 3767.4546
$ws.Range("M88").Value = -1946
$ws.Range("N88").Value = -4579.4546

# Row 91: "The Rose and the Riveter (L)" (Adamantite Rivets)
$ws.Range("H91").Value = 3267.8823
$ws.Range("I91").Value = 2352
$ws.Range("J91").Value = 3767.4546
$ws.Range("K91").Value = 2352
$ws.Range("L91").Value = 3767.4546
$ws.Range("M91").Value = -948
$ws.Range("N91").Value = -6575.4546

# Row 109: "A Head of Demand" (Deepgold Helm of Fending)
$ws.Range("H109").Value = 30795.4
$ws.Range("J109").Value = 30795.4
$ws.Range("L109").Value = 30795.4
$ws.Range("N109").Value = -33569.4

# Row 132: "Don't Bore Me, Ore Me" (Mountain Chromite Ingot)
$ws.Range("H132").Value = 28228.418
$ws.Range("I132").Value = 19691.564
$ws.Range("J132").Value = 47792.043
$ws.Range("K132").Value = 59074.692
$ws.Range("L132").Value = 143376.129
$ws.Range("M132").Value = -56544.692
$ws.Range("N132").Value = -148436.129

$ws = $wb.Worksheets.Item("BSM")
# Row 105: "Ingot to Wing It" (Molybdenum Ingot)
$ws.Range("H105").Value = 71431340
$ws.Range("I105").Value = 71431340
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 71431340
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -71429593
$ws.Range("N105").ClearContents()

# Row 140: "Ceremonial Teeth" (Ra'Kaznar Twinfangs)
$ws.Range("H140").Value = 63740
$ws.Range("J140").Value = 63740
$ws.Range("L140").Value = 63740
$ws.Range("N140").Value = -74100

$ws = $wb.Worksheets.Item("CRP")
# Row 62: "Splinter in the Sewers" (Cedar Lumber)
$ws.Range("H62").Value = 3682
$ws.Range("I62").Value = 3251
$ws.Range("J62").Value = 3928.2856
$ws.Range("K62").Value = 3251
$ws.Range("L62").Value = 3928.2856
$ws.Range("M62").Value = -2627
$ws.Range("N62").Value = -5176.2856

# Row 65: "The Lumber of Their Discontent (L)" (Cedar Lumber)
$ws.Range("H65").Value = 3682
$ws.Range("I65").Value = 3251
$ws.Range("J65").Value = 3928.2856
$ws.Range("K65").Value = 16255
$ws.Range("L65").Value = 19641.428
$ws.Range("M65").Value = -13135
$ws.Range("N65").Value = -25881.428

# Row 99: "O Pine" (Pine Lumber)
$ws.Range("H99").Value = 1725
$ws.Range("I99").Value = 1700
$ws.Range("K99").Value = 1700
$ws.Range("M99").Value = -202

# Row 126: "A Better Conductor" (Red Pine Lumber)
$ws.Range("H126").Value = 1725
$ws.Range("I126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("M126").Value = -2630

$ws = $wb.Worksheets.Item("CUL")
# Row 2: "Pork Is a Salty Food" (Table Salt)
$ws.Range("H2").Value = 531.125
$ws.Range("I2").Value = 578.4286
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 3470.5716
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -3357.5716
$ws.Range("N2").Value = -1426

# Row 131: "The Mountain Steeped" (Tsai tou Vounou)
$ws.Range("H131").Value = 877.9231
$ws.Range("I131").Value = 318.625
$ws.Range("J131").Value = 1022.25806
$ws.Range("K131").Value = 955.875
$ws.Range("L131").Value = 3066.77418
$ws.Range("M131").Value = 4084.125
$ws.Range("N131").Value = -13146.77418

$ws = $wb.Worksheets.Item("GSM")
# Row 70: "Sky Is the Limit" (Mythrite Ingot)
$ws.Range("H70").Value = 61872.223
$ws.Range("I70").Value = 130500
$ws.Range("J70").Value = 6970
$ws.Range("K70").Value = 130500
$ws.Range("L70").Value = 6970
$ws.Range("M70").Value = -130230
$ws.Range("N70").Value = -7510

# Row 73: "Hulls of Broken Dreams (L)" (Mythrite Ingot)
$ws.Range("H73").Value = 61872.223
$ws.Range("I73").Value = 130500
$ws.Range("J73").Value = 6970
$ws.Range("K73").Value = 130500
$ws.Range("L73").Value = 6970
$ws.Range("M73").Value = -129564
$ws.Range("N73").Value = -8842

# Row 80: "Needs More Prayerbell" (Hardsilver Ingot)
$ws.Range("H80").Value = 4418.2144
$ws.Range("I80").Value = 4019
$ws.Range("J80").Value = 4640
$ws.Range("K80").Value = 4019
$ws.Range("L80").Value = 4640
$ws.Range("M80").Value = -3021
$ws.Range("N80").Value = -6636

# Row 83: "With a Noise That Reaches Heaven (L)" (Hardsilver Ingot)
$ws.Range("H83").Value = 4418.2144
$ws.Range("I83").Value = 4019
$ws.Range("J83").Value = 4640
$ws.Range("K83").Value = 20095
$ws.Range("L83").Value = 23200
$ws.Range("M83").Value = -15103
$ws.Range("N83").Value = -33184

# Row 132: "On Board for Lar" (Lar Ingot)
$ws.Range("H132").Value = 73076.96000000001
$ws.Range("I132").Value = 46766.227
$ws.Range("J132").Value = 169549.67
$ws.Range("K132").Value = 140298.681
$ws.Range("L132").Value = 508649.01
$ws.Range("M132").Value = -137768.681
$ws.Range("N132").Value = -513709.01

# Row 136: "Shiny and Good" (Pink Beryl)
$ws.Range("H136").Value = 27064.9
$ws.Range("J136").Value = 27064.9
$ws.Range("L136").Value = 81194.70000000001
$ws.Range("N136").Value = -86294.70000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 128: "Grips of Fear" (Kumbhiraskin Gloves of the Black Griffin)
$ws.Range("H128").Value = 55429
$ws.Range("J128").Value = 55429
$ws.Range("L128").Value = 55429
$ws.Range("N128").Value = -65389

$ws = $wb.Worksheets.Item("WVR")
# Row 137: "Traditional Trousers" (Sarcenet Slops of Aiming)
$ws.Range("H137").Value = 75725
$ws.Range("J137").Value = 75725
$ws.Range("L137").Value = 75725
$ws.Range("N137").Value = -85925
